$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = "4x5090"
$ws.Range("B8").Value = 9790.98
$ws.Range("C8").Value = 2.6
$ws.Range("D8").Value = 0.07376403814758301
